$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 14. This pushes
# all existing data rows (old 14..39) down by two (to 16..41), exactly as
# required: the former last row (old 39) becomes the new row 41.
$ws.Rows(14).Insert()
$ws.Rows(14).Insert()

# Fill in the two newly inserted rows with their data.

# New row 14
$ws.Cells.Item(14, 1).Value = 3
$ws.Cells.Item(14, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 45054
$ws.Cells.Item(14, 5).Value = 5
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100104
$ws.Cells.Item(14, 8).Value = "Frutos de pepita"
$ws.Cells.Item(14, 9).Value = 100104001
$ws.Cells.Item(14, 10).Value = "Granada"
$ws.Cells.Item(14, 11).Value = "Wonderfull"
$ws.Cells.Item(14, 12).Value = "Especial"
$ws.Cells.Item(14, 13).Value = 54
$ws.Cells.Item(14, 14).Value = 16000
$ws.Cells.Item(14, 15).Value = 16000
$ws.Cells.Item(14, 16).Value = 16000
$ws.Cells.Item(14, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(14, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(14, 19).Value = 1143
$ws.Cells.Item(14, 20).Value = 14

# New row 15
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(15, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 45054
$ws.Cells.Item(15, 5).Value = 5
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100104
$ws.Cells.Item(15, 8).Value = "Frutos de pepita"
$ws.Cells.Item(15, 9).Value = 100104001
$ws.Cells.Item(15, 10).Value = "Granada"
$ws.Cells.Item(15, 11).Value = "Wonderfull"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 50
$ws.Cells.Item(15, 14).Value = 14000
$ws.Cells.Item(15, 15).Value = 14000
$ws.Cells.Item(15, 16).Value = 14000
$ws.Cells.Item(15, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(15, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 19).Value = 1000
$ws.Cells.Item(15, 20).Value = 14

# Make sure the new D14/D15 date cells use the same date-time number
# format as the rest of column D.
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("D15").NumberFormat = $ws.Range("D13").NumberFormat
